$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @("trainingaudio/07_pitapi2.wav", "pngimages/07_suitcase.png"),
    @("trainingaudio/15_kopota3.wav", "pngimages/15_barrel.png"),
    @("trainingaudio/21_papika1.wav", "pngimages/21_cheese.png"),
    @("trainingaudio/18_popata2.wav", "pngimages/18_donut.png"),
    @("trainingaudio/01_kitipi1.wav", "pngimages/01_gift.png"),
    @("trainingaudio/06_titoka3.wav", "pngimages/06_tent.png"),
    @("trainingaudio/19_papipi1.wav", "pngimages/19_burger.png"),
    @("trainingaudio/04_kitoti2.wav", "pngimages/04_ladder.png"),
    @("trainingaudio/20_tatito1.wav", "pngimages/20_pizza.png")
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
